$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Norway: duplicate of the "Finland" sheet layout (which has the extra
# MZXSDR240 repeater row), renamed and re-pointed at the Norway market data,
# inserted right after "Hungary".
# ---------------------------------------------------------------------------
$finland = $wb.Worksheets.Item("Finland")
$hungary = $wb.Worksheets.Item("Hungary")

$finland.Copy([System.Reflection.Missing]::Value, $hungary)
$norway = $wb.Worksheets.Item($hungary.Index + 1)
$norway.Name = "Norway"

# Set B4 (user story / NGC code) before B2 (market name) so the shared
# strings table picks up the same ordering as the source edit.
$norway.Range("B4").Value = "NGC-2931/T3071"
$norway.Range("B2").Value = "Norway Market"

$norway.Columns.Item(1).ColumnWidth = 19.65
$norway.Columns.Item(2).ColumnWidth = 15
$norway.Columns.Item(3).ColumnWidth = 12.16
$norway.Columns.Item(4).ColumnWidth = 12.16

$norway.Range("B15").Select()

# ---------------------------------------------------------------------------
# Poland: duplicate of the "Hungary" sheet layout, renamed and re-pointed at
# the Poland market data, inserted right after "Norway".
# ---------------------------------------------------------------------------
$hungary.Copy([System.Reflection.Missing]::Value, $norway)
$poland = $wb.Worksheets.Item($norway.Index + 1)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/3036/T3037"
$poland.Range("B2").Value = "Poland Market"

$poland.Columns.Item(1).ColumnWidth = 19.65
$poland.Columns.Item(2).ColumnWidth = 24
$poland.Columns.Item(3).ColumnWidth = 12.16
$poland.Columns.Item(4).ColumnWidth = 12.16

$poland.Range("B15").Select()

# The source edit leaves "Norway" as the active/selected tab (not the
# last-created "Poland" sheet), so re-activate it last.
$norway.Activate()
